$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new record as the new row 181, pushing the existing
# rows 181:289 down to 182:290 (this matches the diff: dimension grows
# from A1:T289 to A1:T290, and every record from the old row 181 onward
# is now found one row further down, with the former last row (289)
# becoming row 290).
$ws.Rows.Item(181).Insert()

$ws.Range("A181").Value = 6
$ws.Range("B181").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C181").Value = "Metropolitana"
$ws.Range("D181").Value = 44960
$ws.Range("E181").Value = 13
$ws.Range("F181").Value = "Fruta"
$ws.Range("G181").Value = 100101
$ws.Range("H181").Value = "Berries"
$ws.Range("I181").Value = 100101004
$ws.Range("J181").Value = "Frambuesa"
$ws.Range("K181").Value = "Sin especificar"
$ws.Range("L181").Value = "Primera"
$ws.Range("M181").Value = 250
$ws.Range("N181").Value = 7000
$ws.Range("O181").Value = 7000
$ws.Range("P181").Value = 7000
$ws.Range("Q181").Value = '$/bandeja 2 kilos'
$ws.Range("R181").Value = "Región del Maule"
$ws.Range("S181").Value = 3500
$ws.Range("T181").Value = 2
